# Append four new user records (rows 7-10) to the Users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    # Force the cell to be stored as text (not inferred as a number), even
    # when the value looks numeric (phone numbers, zip-like ids, ...).
    # A leading apostrophe makes Excel treat the input as literal text; the
    # following ClearFormats() drops the resulting "quote prefix" style so
    # the cell keeps using the default (General) style, matching every
    # other cell on the sheet.
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Row 7: heyman
Set-TextCell $ws.Cells.Item(7, 1) "heyman"
Set-TextCell $ws.Cells.Item(7, 2) "heyman@gmail.com"
Set-TextCell $ws.Cells.Item(7, 3) "+77254567890"
Set-TextCell $ws.Cells.Item(7, 4) "admin"
Set-TextCell $ws.Cells.Item(7, 5) "101"
Set-TextCell $ws.Cells.Item(7, 6) "Gurgaon"
Set-TextCell $ws.Cells.Item(7, 7) "HR"
Set-TextCell $ws.Cells.Item(7, 8) "IND"

# Row 8: my
Set-TextCell $ws.Cells.Item(8, 1) "my"
Set-TextCell $ws.Cells.Item(8, 2) "my@gmail.com"
Set-TextCell $ws.Cells.Item(8, 3) "+97237777890"
Set-TextCell $ws.Cells.Item(8, 4) "admin"
Set-TextCell $ws.Cells.Item(8, 5) "156"
Set-TextCell $ws.Cells.Item(8, 6) "ghn"
Set-TextCell $ws.Cells.Item(8, 7) "Hill"
Set-TextCell $ws.Cells.Item(8, 8) "Ing"

# Row 9: shu (column E left blank, but touched so an empty cell is recorded)
Set-TextCell $ws.Cells.Item(9, 1) "shu"
Set-TextCell $ws.Cells.Item(9, 2) "shu@gmail.com"
Set-TextCell $ws.Cells.Item(9, 3) "+911234567891"
Set-TextCell $ws.Cells.Item(9, 4) "user"
$ws.Cells.Item(9, 5).ClearFormats()
Set-TextCell $ws.Cells.Item(9, 6) "jai"
Set-TextCell $ws.Cells.Item(9, 7) "RJ"
Set-TextCell $ws.Cells.Item(9, 8) "ind"

# Row 10: doller (column E left blank, but touched so an empty cell is recorded)
Set-TextCell $ws.Cells.Item(10, 1) "doller"
Set-TextCell $ws.Cells.Item(10, 2) "doller@gmail.com"
Set-TextCell $ws.Cells.Item(10, 3) "+88233367890"
Set-TextCell $ws.Cells.Item(10, 4) "user"
$ws.Cells.Item(10, 5).ClearFormats()
Set-TextCell $ws.Cells.Item(10, 6) "bank"
Set-TextCell $ws.Cells.Item(10, 7) "usa"
Set-TextCell $ws.Cells.Item(10, 8) "united"
